$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.04683631701326
$ws.Cells.Item(2, 4).Value = 1.044526583574538
$ws.Cells.Item(2, 5).Value = 1.05404822482355
$ws.Cells.Item(2, 6).Value = 1.064288156663795
$ws.Cells.Item(2, 9).Value = 1.036873849467261
$ws.Cells.Item(2, 10).Value = 1.051888589540881
$ws.Cells.Item(2, 11).Value = 1.047297384801488
$ws.Cells.Item(2, 12).Value = 1.056792509977215
$ws.Cells.Item(2, 13).Value = 1.067004498653303
$ws.Cells.Item(2, 14).Value = 1.053382391991337

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04847530490094
$ws.Cells.Item(3, 4).Value = 1.045723588076897
$ws.Cells.Item(3, 5).Value = 1.055581864906168
$ws.Cells.Item(3, 6).Value = 1.066079459659062
$ws.Cells.Item(3, 9).Value = 1.037261713018638
$ws.Cells.Item(3, 10).Value = 1.053172414506922
$ws.Cells.Item(3, 11).Value = 1.048304639119227
$ws.Cells.Item(3, 12).Value = 1.058137473881384
$ws.Cells.Item(3, 13).Value = 1.06860853398062
$ws.Cells.Item(3, 14).Value = 1.054668040136087

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.049533318452228
$ws.Cells.Item(4, 4).Value = 1.046495689038585
$ws.Cells.Item(4, 5).Value = 1.056572153961129
$ws.Cells.Item(4, 6).Value = 1.067236771621106
$ws.Cells.Item(4, 9).Value = 1.037509883150625
$ws.Cells.Item(4, 10).Value = 1.054000289361745
$ws.Cells.Item(4, 11).Value = 1.048953397801746
$ws.Cells.Item(4, 12).Value = 1.059005183328979
$ws.Cells.Item(4, 13).Value = 1.069644226292876
$ws.Cells.Item(4, 14).Value = 1.055497090668162

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.049977515835262
$ws.Cells.Item(5, 4).Value = 1.046819703636782
$ws.Cells.Item(5, 5).Value = 1.056987984496798
$ws.Cells.Item(5, 6).Value = 1.067722891883113
$ws.Cells.Item(5, 9).Value = 1.037613546011512
$ws.Cells.Item(5, 10).Value = 1.054347656058384
$ws.Cells.Item(5, 11).Value = 1.049225424491637
$ws.Cells.Item(5, 12).Value = 1.059369362433763
$ws.Cells.Item(5, 13).Value = 1.07007911005171
$ws.Cells.Item(5, 14).Value = 1.055844950665354

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.050052064192567
$ws.Cells.Item(6, 4).Value = 1.04687407352201
$ws.Cells.Item(6, 5).Value = 1.057057776024247
$ws.Cells.Item(6, 6).Value = 1.067804489721182
$ws.Cells.Item(6, 9).Value = 1.037630912393535
$ws.Cells.Item(6, 10).Value = 1.054405941326267
$ws.Cells.Item(6, 11).Value = 1.049271057459538
$ws.Cells.Item(6, 12).Value = 1.059430474458896
$ws.Cells.Item(6, 13).Value = 1.070152098714276
$ws.Cells.Item(6, 14).Value = 1.055903318705003

$ws.Cells.Item(7, 2).Value = 1.019999999999999
$ws.Cells.Item(7, 3).Value = 1.04953925614701
$ws.Cells.Item(7, 4).Value = 1.046500020792575
$ws.Cells.Item(7, 5).Value = 1.056577712207903
$ws.Cells.Item(7, 6).Value = 1.067243268791909
$ws.Cells.Item(7, 9).Value = 1.037511270918462
$ws.Cells.Item(7, 10).Value = 1.054004933519644
$ws.Cells.Item(7, 11).Value = 1.048957035421139
$ws.Cells.Item(7, 12).Value = 1.059010051875382
$ws.Cells.Item(7, 13).Value = 1.069650039263297
$ws.Cells.Item(7, 14).Value = 1.055501741421297

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.047390749969239
$ws.Cells.Item(8, 4).Value = 1.044931626950305
$ws.Cells.Item(8, 5).Value = 1.054566962496569
$ws.Cells.Item(8, 6).Value = 1.064893912280936
$ws.Cells.Item(8, 9).Value = 1.037005512259378
$ws.Cells.Item(8, 10).Value = 1.052323059102204
$ws.Cells.Item(8, 11).Value = 1.047638417142532
$ws.Cells.Item(8, 12).Value = 1.057247585341569
$ws.Cells.Item(8, 13).Value = 1.067547057761749
$ws.Cells.Item(8, 14).Value = 1.053817478549271

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.043584954397953
$ws.Cells.Item(9, 4).Value = 1.0421488723906
$ws.Cells.Item(9, 5).Value = 1.051007365176179
$ws.Cells.Item(9, 6).Value = 1.060739836928771
$ws.Cells.Item(9, 9).Value = 1.036092679569106
$ws.Cells.Item(9, 10).Value = 1.049337163287071
$ws.Cells.Item(9, 11).Value = 1.045291521366711
$ws.Cells.Item(9, 12).Value = 1.054121751823649
$ws.Cells.Item(9, 13).Value = 1.063823767444305
$ws.Cells.Item(9, 14).Value = 1.050827342419594

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041033645201152
$ws.Cells.Item(10, 4).Value = 1.040280406201721
$ws.Cells.Item(10, 5).Value = 1.048622600237341
$ws.Cells.Item(10, 6).Value = 1.057960100013911
$ws.Cells.Item(10, 9).Value = 1.035469375643349
$ws.Cells.Item(10, 10).Value = 1.047331032386908
$ws.Cells.Item(10, 11).Value = 1.043710780777549
$ws.Cells.Item(10, 12).Value = 1.052023697773499
$ws.Cells.Item(10, 13).Value = 1.06132902988328
$ws.Cells.Item(10, 14).Value = 1.048818362583446

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.03992538188453
$ws.Cells.Item(11, 4).Value = 1.039468075257443
$ws.Cells.Item(11, 5).Value = 1.047587044639802
$ws.Cells.Item(11, 6).Value = 1.056753808616156
$ws.Cells.Item(11, 9).Value = 1.03519593204117
$ws.Cells.Item(11, 10).Value = 1.046458537617744
$ws.Cells.Item(11, 11).Value = 1.04302236951652
$ws.Cells.Item(11, 12).Value = 1.051111717985073
$ws.Cells.Item(11, 13).Value = 1.060245641846727
$ws.Cells.Item(11, 14).Value = 1.047944628771633

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.039513178423022
$ws.Cells.Item(12, 4).Value = 1.039165838024379
$ws.Cells.Item(12, 5).Value = 1.047201938886261
$ws.Cells.Item(12, 6).Value = 1.0563053251025
$ws.Cells.Item(12, 9).Value = 1.035093825512704
$ws.Cells.Item(12, 10).Value = 1.046133867654765
$ws.Cells.Item(12, 11).Value = 1.042766061870443
$ws.Cells.Item(12, 12).Value = 1.050772428854195
$ws.Cells.Item(12, 13).Value = 1.059842735800597
$ws.Cells.Item(12, 14).Value = 1.047619497740067

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.039601622372523
$ws.Cells.Item(13, 4).Value = 1.039230691780289
$ws.Cells.Item(13, 5).Value = 1.047284566147326
$ws.Cells.Item(13, 6).Value = 1.056401545286366
$ws.Cells.Item(13, 9).Value = 1.035115752096442
$ws.Cells.Item(13, 10).Value = 1.046203537194791
$ws.Cells.Item(13, 11).Value = 1.042821068076504
$ws.Cells.Item(13, 12).Value = 1.050845232121724
$ws.Cells.Item(13, 13).Value = 1.059929182804896
$ws.Cells.Item(13, 14).Value = 1.047689266218831

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.03989132022451
$ws.Cells.Item(14, 4).Value = 1.039443102517163
$ws.Cells.Item(14, 5).Value = 1.047555221021867
$ws.Cells.Item(14, 6).Value = 1.056716745347975
$ws.Cells.Item(14, 9).Value = 1.035187502871067
$ws.Cells.Item(14, 10).Value = 1.046431712330591
$ws.Cells.Item(14, 11).Value = 1.043001195360783
$ws.Cells.Item(14, 12).Value = 1.051083683310371
$ws.Cells.Item(14, 13).Value = 1.060212347534672
$ws.Cells.Item(14, 14).Value = 1.047917765389495

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040069739952005
$ws.Cells.Item(15, 4).Value = 1.039573909057824
$ws.Cells.Item(15, 5).Value = 1.047721919829852
$ws.Cells.Item(15, 6).Value = 1.056910895292209
$ws.Cells.Item(15, 9).Value = 1.035231639553573
$ws.Cells.Item(15, 10).Value = 1.046572220488136
$ws.Cells.Item(15, 11).Value = 1.043112097806602
$ws.Cells.Item(15, 12).Value = 1.051230529139566
$ws.Cells.Item(15, 13).Value = 1.060386749645344
$ws.Cells.Item(15, 14).Value = 1.04805847308474

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.04110712136002
$ws.Cells.Item(16, 4).Value = 1.040334248156281
$ws.Cells.Item(16, 5).Value = 1.048691263665226
$ws.Cells.Item(16, 6).Value = 1.058040100513093
$ws.Cells.Item(16, 9).Value = 1.035487448064768
$ws.Cells.Item(16, 10).Value = 1.047388855375263
$ws.Cells.Item(16, 11).Value = 1.043756384577853
$ws.Cells.Item(16, 12).Value = 1.052084147916643
$ws.Cells.Item(16, 13).Value = 1.061400863198558
$ws.Cells.Item(16, 14).Value = 1.048876267687077

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041756887910492
$ws.Cells.Item(17, 4).Value = 1.040810305793307
$ws.Cells.Item(17, 5).Value = 1.049298512016511
$ws.Cells.Item(17, 6).Value = 1.05874770155262
$ws.Cells.Item(17, 9).Value = 1.035646956885536
$ws.Cells.Item(17, 10).Value = 1.047900076146989
$ws.Cells.Item(17, 11).Value = 1.044159467373231
$ws.Cells.Item(17, 12).Value = 1.052618652576205
$ws.Cells.Item(17, 13).Value = 1.062036136453119
$ws.Cells.Item(17, 14).Value = 1.049388214450936

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042135546201646
$ws.Cells.Item(18, 4).Value = 1.041087667043981
$ws.Cells.Item(18, 5).Value = 1.049652427254739
$ws.Cells.Item(18, 6).Value = 1.059160178928406
$ws.Cells.Item(18, 9).Value = 1.035739653553802
$ws.Cells.Item(18, 10).Value = 1.048197894129881
$ws.Cells.Item(18, 11).Value = 1.044394199200125
$ws.Cells.Item(18, 12).Value = 1.052930082402836
$ws.Cells.Item(18, 13).Value = 1.062406377612387
$ws.Cells.Item(18, 14).Value = 1.049686455369525

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042264601747298
$ws.Cells.Item(19, 4).Value = 1.041182186904121
$ws.Cells.Item(19, 5).Value = 1.049773055690443
$ws.Cells.Item(19, 6).Value = 1.059300780380351
$ws.Cells.Item(19, 9).Value = 1.035771202817992
$ws.Cells.Item(19, 10).Value = 1.048299380220607
$ws.Cells.Item(19, 11).Value = 1.044474172582338
$ws.Cells.Item(19, 12).Value = 1.053036215010583
$ws.Cells.Item(19, 13).Value = 1.062532569316294
$ws.Cells.Item(19, 14).Value = 1.04978808558214

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041687209336827
$ws.Cells.Item(20, 4).Value = 1.040759261977795
$ws.Cells.Item(20, 5).Value = 1.049233389325265
$ws.Cells.Item(20, 6).Value = 1.058671809077204
$ws.Cells.Item(20, 9).Value = 1.035629878517049
$ws.Cells.Item(20, 10).Value = 1.047845265190152
$ws.Cells.Item(20, 11).Value = 1.044116259728119
$ws.Cells.Item(20, 12).Value = 1.05256134030196
$ws.Cells.Item(20, 13).Value = 1.061968009086585
$ws.Cells.Item(20, 14).Value = 1.049333325656254

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.039806026595564
$ws.Cells.Item(21, 4).Value = 1.039380566767218
$ws.Cells.Item(21, 5).Value = 1.047475532553902
$ws.Cells.Item(21, 6).Value = 1.056623938317015
$ws.Cells.Item(21, 9).Value = 1.03516638892529
$ws.Cells.Item(21, 10).Value = 1.04636453672601
$ws.Cells.Item(21, 11).Value = 1.042948169030615
$ws.Cells.Item(21, 12).Value = 1.051013480361902
$ws.Cells.Item(21, 13).Value = 1.060128976151438
$ws.Cells.Item(21, 14).Value = 1.04785049438785

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03862009097946
$ws.Cells.Item(22, 4).Value = 1.038510820838061
$ws.Cells.Item(22, 5).Value = 1.046367663223887
$ws.Cells.Item(22, 6).Value = 1.055333962372633
$ws.Cells.Item(22, 9).Value = 1.034871863265579
$ws.Cells.Item(22, 10).Value = 1.045430144480179
$ws.Cells.Item(22, 11).Value = 1.042210261965814
$ws.Cells.Item(22, 12).Value = 1.050037153461889
$ws.Cells.Item(22, 13).Value = 1.058969875532991
$ws.Cells.Item(22, 14).Value = 1.046914775197855

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.039249082543405
$ws.Cells.Item(23, 4).Value = 1.038972168275034
$ws.Cells.Item(23, 5).Value = 1.0469552198946
$ws.Cells.Item(23, 6).Value = 1.05601803557716
$ws.Cells.Item(23, 9).Value = 1.035028293233686
$ws.Cells.Item(23, 10).Value = 1.045925809632727
$ws.Cells.Item(23, 11).Value = 1.0426017735852
$ws.Cells.Item(23, 12).Value = 1.050555023202746
$ws.Cells.Item(23, 13).Value = 1.059584609509203
$ws.Cells.Item(23, 14).Value = 1.047411144251773

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.04171869513408
$ws.Cells.Item(24, 4).Value = 1.040782327454614
$ws.Cells.Item(24, 5).Value = 1.049262816338997
$ws.Cells.Item(24, 6).Value = 1.058706102405909
$ws.Cells.Item(24, 9).Value = 1.035637596553818
$ws.Cells.Item(24, 10).Value = 1.047870033040683
$ws.Cells.Item(24, 11).Value = 1.044135784575892
$ws.Cells.Item(24, 12).Value = 1.052587238295143
$ws.Cells.Item(24, 13).Value = 1.061998793847147
$ws.Cells.Item(24, 14).Value = 1.049358128679974

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.044571274661826
$ws.Cells.Item(25, 4).Value = 1.042870588311257
$ws.Cells.Item(25, 5).Value = 1.05192961879905
$ws.Cells.Item(25, 6).Value = 1.0618155307087
$ws.Cells.Item(25, 9).Value = 1.036331252573766
$ws.Cells.Item(25, 10).Value = 1.050111782808785
$ws.Cells.Item(25, 11).Value = 1.045901061291052
$ws.Cells.Item(25, 12).Value = 1.054932305379245
$ws.Cells.Item(25, 13).Value = 1.064788484095818
$ws.Cells.Item(25, 14).Value = 1.051603061989879
